$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update destination labels (rows 8-20 shift up by one; row 8 is new, row 21 is appended) ---
$ws.Range("A8").Value = "Tracy, CA 95304"
$ws.Range("A9").Value = "Fort Worth, TX 76140"
$ws.Range("A10").Value = "Modesto, CA 95353"
$ws.Range("A11").Value = "Dallas, TX 75244"
$ws.Range("A12").Value = "Chicago, IL 60628"
$ws.Range("A13").Value = "Stockton, CA 95205"
$ws.Range("A14").Value = "Stockton, CA 96215"
$ws.Range("A15").Value = "Fremont, CA 94538"
$ws.Range("A16").Value = "Pleasant Grove, CA 95668"
$ws.Range("A17").Value = "Fairfield, CA 94533"
$ws.Range("A18").Value = "Manteca, CA 95336"
$ws.Range("A19").Value = "Milpitas, CA 95035"
$ws.Range("A20").Value = "Boise, ID 83717"

# --- Add new row 21 (Loveland, CO 80538) ---
$ws.Range("A21").Value = "Loveland, CO 80538"
$ws.Range("B21").NumberFormat = $ws.Range("B20").NumberFormat

# --- Update numeric rate values ---
$ws.Range("B2").Value = 440
$ws.Range("C2").Value = 830
$ws.Range("D2").Value = 1190

$ws.Range("B3").Value = 425
$ws.Range("C3").Value = 805
$ws.Range("D3").Value = 1340

$ws.Range("B4").Value = 450
$ws.Range("C4").Value = 840
$ws.Range("D4").Value = 1190

$ws.Range("B5").Value = 515
$ws.Range("C5").Value = 955
$ws.Range("D5").Value = 1315

$ws.Range("B6").Value = 570
$ws.Range("C6").Value = 1085
$ws.Range("D6").Value = 1535

$ws.Range("B7").Value = 655
$ws.Range("C7").Value = 1245
$ws.Range("D7").Value = 1790

$ws.Range("B8").Value = 395
$ws.Range("C8").Value = 530
$ws.Range("D8").Value = 665
$ws.Range("E8").Value = 785
$ws.Range("F8").Value = 880
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("J8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("O8").ClearContents()
$ws.Range("P8").ClearContents()

$ws.Range("B9").Value = 485
$ws.Range("C9").ClearContents()
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()

$ws.Range("B10").Value = 395
$ws.Range("C10").Value = 530
$ws.Range("D10").Value = 665
$ws.Range("E10").Value = 785
$ws.Range("F10").Value = 880
$ws.Range("G10").Value = 980

$ws.Range("B11").ClearContents()
$ws.Range("C11").ClearContents()
$ws.Range("D11").ClearContents()
$ws.Range("E11").ClearContents()
$ws.Range("F11").ClearContents()
$ws.Range("G11").Value = 1800

$ws.Range("E12").Value = 1435
$ws.Range("G12").ClearContents()

$ws.Range("B13").Value = 335
$ws.Range("C13").Value = 480
$ws.Range("D13").Value = 610
$ws.Range("E13").Value = 700
$ws.Range("F13").Value = 840
$ws.Range("G13").Value = 975

$ws.Range("B14").Value = 400
$ws.Range("C14").Value = 540
$ws.Range("D14").Value = 690
$ws.Range("E14").Value = 820
$ws.Range("F14").Value = 910
$ws.Range("G14").Value = 1000
$ws.Range("H14").Value = 1050
$ws.Range("I14").Value = 1105
$ws.Range("J14").Value = 1150
$ws.Range("K14").Value = 1200

$ws.Range("B15").Value = 395
$ws.Range("C15").Value = 530
$ws.Range("D15").Value = 665
$ws.Range("E15").Value = 785
$ws.Range("F15").Value = 880
$ws.Range("G15").Value = 980
$ws.Range("H15").ClearContents()
$ws.Range("I15").ClearContents()
$ws.Range("J15").ClearContents()
$ws.Range("K15").ClearContents()

$ws.Range("B16").Value = 395
$ws.Range("C16").Value = 530
$ws.Range("D16").Value = 665
$ws.Range("E16").Value = 785
$ws.Range("F16").Value = 880
$ws.Range("G16").Value = 980

$ws.Range("B17").Value = 400
$ws.Range("C17").Value = 565
$ws.Range("D17").Value = 705
$ws.Range("E17").Value = 835
$ws.Range("F17").Value = 930
$ws.Range("G17").Value = 1030

$ws.Range("B18").Value = 400
$ws.Range("C18").Value = 565
$ws.Range("D18").Value = 705
$ws.Range("E18").Value = 835
$ws.Range("F18").Value = 930
$ws.Range("G18").Value = 1030

$ws.Range("B19").Value = 400
$ws.Range("C19").Value = 565
$ws.Range("D19").Value = 705
$ws.Range("E19").Value = 835
$ws.Range("F19").Value = 930
$ws.Range("G19").Value = 1030

$ws.Range("B20").ClearContents()
$ws.Range("C20").Value = 1195
$ws.Range("D20").ClearContents()
$ws.Range("E20").ClearContents()
$ws.Range("F20").ClearContents()
$ws.Range("G20").ClearContents()

# --- Row 21 value ---
$ws.Range("B21").Value = 540

# --- Set active selection to K7 (matches final workbook state) ---
$ws.Range("K7").Select()
